$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (F1:H1) with the same style as the existing header row
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header formatting (bold font, border, centered alignment) from an
# existing header cell onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Fill in the new "Outliers_MAD" boolean columns for all data rows with FALSE
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("F$r").Value = $false
    $ws.Range("G$r").Value = $false
    $ws.Range("H$r").Value = $false
}
